$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(96, 8).Value = 2565635.8
$ws.Cells.Item(96, 9).Value = 3663932.5
$ws.Cells.Item(96, 10).Value = 2943
$ws.Cells.Item(96, 11).Value = 10991797.5
$ws.Cells.Item(96, 12).Value = 8829
$ws.Cells.Item(96, 13).Value = -10990424.5
$ws.Cells.Item(96, 14).Value = -11575
$ws.Cells.Item(100, 8).Value = 2227.1667
$ws.Cells.Item(100, 9).Value = 1789.909
$ws.Cells.Item(100, 10).Value = 2914.2856
$ws.Cells.Item(100, 11).Value = 1789.909
$ws.Cells.Item(100, 12).Value = 2914.2856
$ws.Cells.Item(100, 13).Value = -1248.909
$ws.Cells.Item(100, 14).Value = -3996.2856
$ws.Cells.Item(137, 8).Value = 1802.38
$ws.Cells.Item(137, 9).Value = 1968
$ws.Cells.Item(137, 10).Value = 1622.9584
$ws.Cells.Item(137, 11).Value = 5904
$ws.Cells.Item(137, 12).Value = 4868.8752
$ws.Cells.Item(137, 13).Value = -3354
$ws.Cells.Item(137, 14).Value = -9968.8752

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6808.3647
$ws.Cells.Item(32, 9).Value = 7184.1714
$ws.Cells.Item(32, 11).Value = 7184.1714
$ws.Cells.Item(32, 13).Value = -6897.1714
$ws.Cells.Item(102, 8).Value = 3316.6667
$ws.Cells.Item(102, 9).Value = 5000
$ws.Cells.Item(102, 11).Value = 5000
$ws.Cells.Item(102, 13).Value = -3378
$ws.Cells.Item(110, 8).Value = 2813.6365
$ws.Cells.Item(110, 9).Value = 2612.5
$ws.Cells.Item(110, 10).Value = 2928.5715
$ws.Cells.Item(110, 11).Value = 2612.5
$ws.Cells.Item(110, 12).Value = 2928.5715
$ws.Cells.Item(110, 13).Value = -567.5
$ws.Cells.Item(110, 14).Value = -7018.5715
$ws.Cells.Item(122, 8).Value = 1188.3636
$ws.Cells.Item(122, 9).Value = 963.55554
$ws.Cells.Item(122, 10).Value = 2200
$ws.Cells.Item(122, 11).Value = 2890.66662
$ws.Cells.Item(122, 12).Value = 6600
$ws.Cells.Item(122, 13).Value = -440.66662
$ws.Cells.Item(122, 14).Value = -11500
$ws.Cells.Item(124, 8).Value = 20602.637
$ws.Cells.Item(124, 10).Value = 20602.637
$ws.Cells.Item(124, 12).Value = 20602.637
$ws.Cells.Item(124, 14).Value = -30422.637
$ws.Cells.Item(125, 8).Value = 39333.332
$ws.Cells.Item(125, 10).Value = 39333.332
$ws.Cells.Item(125, 12).Value = 39333.332
$ws.Cells.Item(125, 14).Value = -49173.332
$ws.Cells.Item(132, 8).Value = 3914.5874
$ws.Cells.Item(132, 9).Value = 2410.7942
$ws.Cells.Item(132, 10).Value = 5677.6553
$ws.Cells.Item(132, 11).Value = 7232.382599999999
$ws.Cells.Item(132, 12).Value = 17032.9659
$ws.Cells.Item(132, 13).Value = -4702.382599999999
$ws.Cells.Item(132, 14).Value = -22092.9659

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 2153.4119
$ws.Cells.Item(107, 9).Value = 2199.6155
$ws.Cells.Item(107, 10).Value = 2003.25
$ws.Cells.Item(107, 11).Value = 2199.6155
$ws.Cells.Item(107, 12).Value = 2003.25
$ws.Cells.Item(107, 13).Value = -279.6154999999999
$ws.Cells.Item(107, 14).Value = -5843.25
$ws.Cells.Item(134, 8).Value = 5887.5835
$ws.Cells.Item(134, 9).Value = 2501.7778
$ws.Cells.Item(134, 10).Value = 9273.388999999999
$ws.Cells.Item(134, 11).Value = 7505.3334
$ws.Cells.Item(134, 12).Value = 27820.167
$ws.Cells.Item(134, 13).Value = -4970.3334
$ws.Cells.Item(134, 14).Value = -32890.167

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 6182.2
$ws.Cells.Item(16, 9).Value = 7637
$ws.Cells.Item(16, 10).Value = 4000
$ws.Cells.Item(16, 11).Value = 7637
$ws.Cells.Item(16, 12).Value = 4000
$ws.Cells.Item(16, 13).Value = -7350
$ws.Cells.Item(16, 14).Value = -4574
$ws.Cells.Item(31, 8).Value = 8132124.5
$ws.Cells.Item(31, 9).Value = 1474.0938
$ws.Cells.Item(31, 10).Value = 37041104
$ws.Cells.Item(31, 11).Value = 1474.0938
$ws.Cells.Item(31, 12).Value = 37041104
$ws.Cells.Item(31, 13).Value = -1179.0938
$ws.Cells.Item(31, 14).Value = -37041694
$ws.Cells.Item(34, 8).Value = 8132124.5
$ws.Cells.Item(34, 9).Value = 1474.0938
$ws.Cells.Item(34, 10).Value = 37041104
$ws.Cells.Item(34, 11).Value = 1474.0938
$ws.Cells.Item(34, 12).Value = 37041104
$ws.Cells.Item(34, 13).Value = -1272.0938
$ws.Cells.Item(34, 14).Value = -37041508
$ws.Cells.Item(105, 8).Value = 2080
$ws.Cells.Item(105, 9).Value = 1904.5454
$ws.Cells.Item(105, 11).Value = 1904.5454
$ws.Cells.Item(105, 13).Value = -157.5454
$ws.Cells.Item(107, 8).Value = 882.64
$ws.Cells.Item(107, 9).Value = 298.92856
$ws.Cells.Item(107, 10).Value = 1625.5454
$ws.Cells.Item(107, 11).Value = 298.92856
$ws.Cells.Item(107, 12).Value = 1625.5454
$ws.Cells.Item(107, 13).Value = 1621.07144
$ws.Cells.Item(107, 14).Value = -5465.5454
$ws.Cells.Item(112, 8).Value = 22856.521
$ws.Cells.Item(112, 10).Value = 22856.521
$ws.Cells.Item(112, 12).Value = 22856.521
$ws.Cells.Item(112, 14).Value = -25810.521
$ws.Cells.Item(113, 8).Value = 6182.2
$ws.Cells.Item(113, 9).Value = 7637
$ws.Cells.Item(113, 10).Value = 4000
$ws.Cells.Item(113, 11).Value = 7637
$ws.Cells.Item(113, 12).Value = 4000
$ws.Cells.Item(113, 13).Value = -5467
$ws.Cells.Item(113, 14).Value = -8340

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(38, 8).Value = 128.125
$ws.Cells.Item(38, 10).Value = 170
$ws.Cells.Item(38, 12).Value = 510
$ws.Cells.Item(38, 14).Value = -1204
$ws.Cells.Item(70, 8).Value = 2243.8333
$ws.Cells.Item(70, 9).Value = 1252
$ws.Cells.Item(70, 10).Value = 3235.6667
$ws.Cells.Item(70, 11).Value = 3756
$ws.Cells.Item(70, 12).Value = 9707.000100000001
$ws.Cells.Item(70, 13).Value = -3441
$ws.Cells.Item(70, 14).Value = -10337.0001
$ws.Cells.Item(73, 8).Value = 2243.8333
$ws.Cells.Item(73, 9).Value = 1252
$ws.Cells.Item(73, 10).Value = 3235.6667
$ws.Cells.Item(73, 11).Value = 3756
$ws.Cells.Item(73, 12).Value = 9707.000100000001
$ws.Cells.Item(73, 13).Value = -2664
$ws.Cells.Item(73, 14).Value = -11891.0001
$ws.Cells.Item(76, 8).Value = 3346
$ws.Cells.Item(76, 9).Value = 2653.3333
$ws.Cells.Item(76, 10).Value = 3642.8572
$ws.Cells.Item(76, 11).Value = 7959.999899999999
$ws.Cells.Item(76, 12).Value = 10928.5716
$ws.Cells.Item(76, 13).Value = -7576.999899999999
$ws.Cells.Item(76, 14).Value = -11694.5716
$ws.Cells.Item(79, 8).Value = 3346
$ws.Cells.Item(79, 9).Value = 2653.3333
$ws.Cells.Item(79, 10).Value = 3642.8572
$ws.Cells.Item(79, 11).Value = 7959.999899999999
$ws.Cells.Item(79, 12).Value = 10928.5716
$ws.Cells.Item(79, 13).Value = -6633.999899999999
$ws.Cells.Item(79, 14).Value = -13580.5716
$ws.Cells.Item(107, 8).Value = 954.5454999999999
$ws.Cells.Item(107, 9).Value = 500
$ws.Cells.Item(107, 11).Value = 1500
$ws.Cells.Item(107, 13).Value = 420
$ws.Cells.Item(131, 8).Value = 1292.7675
$ws.Cells.Item(131, 9).Value = 698.8889
$ws.Cells.Item(131, 10).Value = 1449.9706
$ws.Cells.Item(131, 11).Value = 2096.6667
$ws.Cells.Item(131, 12).Value = 4349.9118
$ws.Cells.Item(131, 13).Value = 2943.3333
$ws.Cells.Item(131, 14).Value = -14429.9118

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 4447.222
$ws.Cells.Item(80, 9).Value = 2754.1667
$ws.Cells.Item(80, 11).Value = 2754.1667
$ws.Cells.Item(80, 13).Value = -1756.1667
$ws.Cells.Item(83, 8).Value = 4447.222
$ws.Cells.Item(83, 9).Value = 2754.1667
$ws.Cells.Item(83, 11).Value = 13770.8335
$ws.Cells.Item(83, 13).Value = -8778.833500000001
$ws.Cells.Item(113, 8).Value = 2697
$ws.Cells.Item(113, 9).Value = 2525
$ws.Cells.Item(113, 10).Value = 2926.3333
$ws.Cells.Item(113, 11).Value = 2525
$ws.Cells.Item(113, 12).Value = 2926.3333
$ws.Cells.Item(113, 13).Value = -355
$ws.Cells.Item(113, 14).Value = -7266.3333
$ws.Cells.Item(122, 8).Value = 2353.3914
$ws.Cells.Item(122, 9).Value = 2116
$ws.Cells.Item(122, 10).Value = 3936
$ws.Cells.Item(122, 11).Value = 6348
$ws.Cells.Item(122, 12).Value = 11808
$ws.Cells.Item(122, 13).Value = -3898
$ws.Cells.Item(122, 14).Value = -16708
$ws.Cells.Item(126, 8).Value = 2885.6072
$ws.Cells.Item(126, 9).Value = 3124.2
$ws.Cells.Item(126, 10).Value = 2289.125
$ws.Cells.Item(126, 11).Value = 9372.599999999999
$ws.Cells.Item(126, 12).Value = 6867.375
$ws.Cells.Item(126, 13).Value = -6902.599999999999
$ws.Cells.Item(126, 14).Value = -11807.375
$ws.Cells.Item(132, 8).Value = 1986738.8
$ws.Cells.Item(132, 9).Value = 2605663.5
$ws.Cells.Item(132, 10).Value = 6179.6
$ws.Cells.Item(132, 11).Value = 7816990.5
$ws.Cells.Item(132, 12).Value = 18538.8
$ws.Cells.Item(132, 13).Value = -7814460.5
$ws.Cells.Item(132, 14).Value = -23598.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(70, 8).Value = 12000
$ws.Cells.Item(70, 10).Value = 15000
$ws.Cells.Item(70, 12).Value = 15000
$ws.Cells.Item(70, 14).Value = -15630
$ws.Cells.Item(73, 8).Value = 12000
$ws.Cells.Item(73, 10).Value = 15000
$ws.Cells.Item(73, 12).Value = 15000
$ws.Cells.Item(73, 14).Value = -17184
$ws.Cells.Item(100, 8).Value = 918.0454999999999
$ws.Cells.Item(100, 9).Value = 846.2857
$ws.Cells.Item(100, 10).Value = 1043.625
$ws.Cells.Item(100, 11).Value = 1692.5714
$ws.Cells.Item(100, 12).Value = 2087.25
$ws.Cells.Item(100, 13).Value = -1151.5714
$ws.Cells.Item(100, 14).Value = -3169.25
$ws.Cells.Item(122, 8).Value = 3182.2273
$ws.Cells.Item(122, 9).Value = 2705.889
$ws.Cells.Item(122, 10).Value = 3938.7646
$ws.Cells.Item(122, 11).Value = 8117.667
$ws.Cells.Item(122, 12).Value = 11816.2938
$ws.Cells.Item(122, 13).Value = -5667.667
$ws.Cells.Item(122, 14).Value = -16716.2938
